$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet was shared/unprotected so collaborators could edit it directly.
[void]$wb.Unprotect()

# --- Row 1: refresh the timestamp, drop the now-unused duplicate in P1 ---
$ws.Range("K1").Value = 43147.51219907407
$ws.Range("P1").ClearContents()

# --- Row 2: a freshly pasted Instagram link + the moment it was added ---
$row2Text = "`nokul (@tlgkyck) " + [char]0x2022 + " Instagram photos and videos`n"
$ws.Range("A2").Value = $row2Text
$ws.Range("K2").Value = 43147.56716963602
$ws.Range("K2").NumberFormat = "yyyy-mm-dd h:mm:ss"
[void]$ws.Rows.Item(2).AutoFit()

# --- Row 3: the original link pasted again, with its own timestamp ---
$ws.Range("A3").Value = $ws.Range("A1").Value2
$ws.Range("K3").Value = 43147.56776974755
$ws.Range("K3").NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
[void]$ws.Rows.Item(3).AutoFit()

# Column K holds the timestamps, so size it to fit its new contents.
[void]$ws.Columns.Item(11).AutoFit()

# Leave the cursor where the last edit happened.
[void]$ws.Range("Q3").Select()
